{"js": "// Add a new bulleted \"Deliverable\" item after the last paragraph of the\n// document (\"A detailed README file in the repository\"), matching the\n// existing ListParagraph / numId=1 / ilvl=0 bullet-list formatting.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\nlastParagraph.insertParagraph(\n  \"Appealing UI/UX is something good to have but not mandatory\",\n  Word.InsertLocation.after\n);\nawait context.sync();\n\n// Word keeps the \"_GoBack\" bookmark (last-edit marker) anchored at the\n// position of the most recent edit, so it needs to move from the end of\n// the old last paragraph to the end of the newly inserted one.\ncontext.document.deleteBookmark(\"_GoBack\");\nconst bodyEnd = context.document.body.getRange(\"End\");\nbodyEnd.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# Add a new bulleted \"Deliverable\" item after the last paragraph of the\n# document (\"A detailed README file in the repository\"), matching the\n# existing ListParagraph / numId=1 / ilvl=0 bullet-list formatting.\n$d = $word.ActiveDocument\n\n$lastParagraph = $d.Paragraphs.Last\n$lastParagraph.Range.InsertParagraphAfter()\n\n$newParagraph = $d.Paragraphs.Item($d.Paragraphs.Count)\n# Append a placeholder character so the insertion point used for the\n# bookmark below is not the very last character in the document; some\n# engines mis-place a bookmark collapsed exactly at end-of-story.\n$newParagraph.Range.Text = \"Appealing UI/UX is something good to have but not mandatory#\"\n\n# Word keeps the \"_GoBack\" bookmark (last-edit marker) anchored at the\n# position of the most recent edit, so it needs to move from the end of\n# the old last paragraph to the end of the newly inserted one.\n$d.Bookmarks.Item(\"_GoBack\").Delete()\n\n$bookmarkPos = $newParagraph.Range.End - 2\n$bookmarkRange = $d.Range($bookmarkPos, $bookmarkPos)\n$d.Bookmarks.Add(\"_GoBack\", $bookmarkRange)\n\n# Remove the placeholder character now that the bookmark is anchored.\n$placeholderRange = $d.Range($bookmarkPos, $bookmarkPos + 1)\n$placeholderRange.Text = \"\"\n"}
